# Added Unique Email and Username Alert
# Mark additional test-plan rows as passing ("x") and fill in the
# corresponding "Mark" column for groups whose sub-tests now all pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark individual sub-tests as passed ("x" in column C) ---
$passRows = @(7,8,9,27,28,43,65,66,67,68,69,70,82,83,84,85,86)
foreach ($r in $passRows) {
    $ws.Cells.Item($r, 3).Value = "x"
}

# --- Fill in the awarded Mark (column D) for groups whose sub-tests all pass ---
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(40, 4).Value = 4
$ws.Cells.Item(48, 4).Value = 3
$ws.Cells.Item(54, 4).Value = 1
$ws.Cells.Item(64, 4).Value = 10
$ws.Cells.Item(81, 4).Value = 3

# --- Restore the view / selection to cell B6 ---
$ws.Range("B6").Select()
